# Commit: "Added support for longer quotes, fixed surplus numnber"
#
# Core content change: the "Surcharge" column (K) on the "Customer Quote"
# sheet used a hard-coded stainless-steel surcharge factor (1.0565) on
# several line items; those were corrected down to 1 (no surcharge),
# matching the "Others @ 1.0" option described in the sheet's own
# "Source options" legend.
#
# There is also a view-state change recorded in the workbook: the
# previously active cell/selection (A28) is moved to G7, and the window
# had scrolled so column B is the left-most visible column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")
$ws.Activate()

# --- Fix the surplus/surcharge factor -------------------------------------
# Rows 16, 17, 20, 23, 24 and 28 had the stainless-steel surcharge
# (1.0565) left in by mistake; reset them to 1 like the rest of the table.
$surchargeRows = 16, 17, 20, 23, 24, 28
foreach ($row in $surchargeRows) {
    $ws.Range("K$row").Value = 1
}

# --- Update the saved view/selection state --------------------------------
# Scroll the window so column B is the left-most visible column (topLeftCell
# = B1) and move the selection to G7.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("G7").Select()
